# Updated cryptos list on Sun May 12 08:25:26 UTC 2024 with GitHub Actions
# Refresh Price (D) and Volume(1h) (E) text columns to the latest scrape values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.837.27"
$ws.Range("E2").Value = "  -0.08%  "
$ws.Range("D3").Value = "2.905.92"
$ws.Range("E3").Value = "  -0.25%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'588.23"
$ws.Range("E5").Value = "  +0.41%  "
$ws.Range("D6").Value = "'144.29"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "'0.504"
$ws.Range("E8").Value = "  -0.15%  "
$ws.Range("D9").Value = "'6.88"
$ws.Range("E9").Value = "  +3.41%  "
$ws.Range("D10").Value = "'0.140"
$ws.Range("E10").Value = "  -2.57%  "
$ws.Range("D11").Value = "'0.438"
$ws.Range("E11").Value = "  -2.19%  "
$ws.Range("E12").Value = "  -0.80%  "
$ws.Range("D13").Value = "'33.31"
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("E14").Value = "  +0.03%  "
$ws.Range("D15").Value = "3.387.28"
$ws.Range("E15").Value = "  -0.23%  "
$ws.Range("D16").Value = "60.776.26"
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "'6.66"
$ws.Range("E17").Value = "  -1.17%  "
$ws.Range("D18").Value = "2.906.07"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "'431.37"
$ws.Range("E19").Value = "  +0.75%  "
$ws.Range("D20").Value = "'13.31"
$ws.Range("E20").Value = "  -1.74%  "
$ws.Range("D21").Value = "'0.674"
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("D22").Value = "'7.08"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "'81.46"
$ws.Range("E23").Value = "  +1.01%  "
$ws.Range("D24").Value = "'10.79"
$ws.Range("E24").Value = "  +1.01%  "
$ws.Range("D25").Value = "'2.17"
$ws.Range("E25").Value = "  -2.40%  "
$ws.Range("D26").Value = "'11.74"
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  +5.08%  "
$ws.Range("D29").Value = "'2.58"
$ws.Range("E29").Value = "  -0.97%  "
$ws.Range("D30").Value = "'6.93"
$ws.Range("E30").Value = "  -4.22%  "
$ws.Range("D31").Value = "'26.42"
$ws.Range("E31").Value = "  -0.61%  "
$ws.Range("D32").Value = "'0.108"
$ws.Range("E32").Value = "  +1.57%  "
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").Value = "0.0₃0861"
$ws.Range("E34").Value = "  -1.34%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  -1.13%  "
$ws.Range("E36").Value = "  -0.24%  "
$ws.Range("D37").Value = "'2.98"
$ws.Range("E37").Value = "  -0.64%  "
$ws.Range("D38").Value = "'1.96"
$ws.Range("E38").Value = "  -1.40%  "
$ws.Range("E39").Value = "  -4.18%  "
$ws.Range("D40").Value = "'8.53"
$ws.Range("E40").Value = "  -0.83%  "
$ws.Range("D41").Value = "'40.96"
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("D42").Value = "'0.280"
$ws.Range("E42").Value = "  -5.79%  "
$ws.Range("D43").Value = "'375.91"
$ws.Range("E43").Value = "  -0.76%  "
$ws.Range("E44").Value = "  -1.85%  "
$ws.Range("D45").Value = "2.693.85"
$ws.Range("E45").Value = "  +0.31%  "
$ws.Range("D46").Value = "'133.43"
$ws.Range("E46").Value = "  +0.44%  "
$ws.Range("D48").Value = "'23.59"
$ws.Range("E48").Value = "  -3.66%  "
$ws.Range("D49").Value = "'0.105"
$ws.Range("E49").Value = "  -0.81%  "
$ws.Range("D50").Value = "'1.99"
$ws.Range("E50").Value = "  -2.82%  "
$ws.Range("E51").Value = "  -1.05%  "
